$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-12 Thursday", "2026-02-13 Friday"),
    @("27×79=", "91×35="),
    @("56×38=", "22×33="),
    @("28×84=", "39×20="),
    @("65×59=", "30×65="),
    @("87×32=", "62×23="),
    @("11×98=", "38×95="),
    @("99×26=", "78×91="),
    @("24×49=", "13×36="),
    @("28×68=", "93×16="),
    @("14×91=", "62×31="),
    @("29×64=", "29×88="),
    @("67×53=", "29×27="),
    @("30×16=", "85×54="),
    @("55×11=", "16×73="),
    @("68×83=", "88×41="),
    @("24×11=", "81×73="),
    @("98×64=", "13×39="),
    @("11×89=", "90×14="),
    @("52×16=", "80×61="),
    @("56×23=", "21×36="),
    @("72×74=", "44×66="),
    @("70×34=", "74×41="),
    @("19×57=", "12×25="),
    @("28×37=", "65×75="),
    @("77×47=", "14×14=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
